$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("H2").Value = [double]"0.001206652733196978"
$ws.Range("I2").Value = [double]"0.001206652733196978"
$ws.Range("L2").Value = [double]"43.02795316867243"
$ws.Range("M2").Value = "[16.639945770677855, 69.41596056666701]"
$ws.Range("N2").Value = [double]"0.001984689032197728"
$ws.Range("O2").Value = [double]"0.001984689032197728"
$ws.Range("Q2").Value = "[0.8365001460008861, 2.4466056901830404]"
$ws.Range("R2").Value = [double]"0.0001669965727382206"
$ws.Range("S2").Value = [double]"0.0001669965727382206"
$ws.Range("T2").Value = [double]"62.44596375843581"
$ws.Range("U2").Value = "[46.6645555150225, 78.22737200184912]"
$ws.Range("V2").Value = [double]"3.815956439723323e-10"
$ws.Range("W2").Value = [double]"3.815956439723323e-10"
$ws.Range("Y2").Value = [double]"15.20420420420449"
$ws.Range("Z2").Value = [double]"21.58498498498538"

# Row 3
$ws.Range("H3").Value = [double]"0.01787431543618312"
$ws.Range("I3").Value = [double]"0.01787431543618312"
$ws.Range("L3").Value = [double]"32.30812269821571"
$ws.Range("M3").Value = "[4.7199221892777885, 59.896323207153635]"
$ws.Range("N3").Value = [double]"0.0227426240641817"
$ws.Range("O3").Value = [double]"0.0227426240641817"
$ws.Range("P3").Value = [double]"1.377394977249579"
$ws.Range("Q3").Value = "[0.20755266780473125, 2.5472372866944264]"
$ws.Range("R3").Value = [double]"0.02205865176989619"
$ws.Range("S3").Value = [double]"0.02205865176989619"
$ws.Range("T3").Value = [double]"68.67262797800859"
$ws.Range("U3").Value = "[52.95620239432, 84.38905356169718]"
$ws.Range("V3").Value = [double]"2.428568457446545e-11"
$ws.Range("W3").Value = [double]"2.428568457446545e-11"
$ws.Range("X3").Value = [double]"19.4414414414418"
$ws.Range("Y3").Value = [double]"14.80540540540568"
$ws.Range("Z3").Value = [double]"24.07747747747791"

# Row 4
$ws.Range("H4").Value = [double]"4.094260309339681e-06"
$ws.Range("I4").Value = [double]"4.094260309339681e-06"
$ws.Range("L4").Value = [double]"49.98008519543264"
$ws.Range("M4").Value = "[26.685333317699047, 73.27483707316624]"
$ws.Range("N4").Value = [double]"8.459728460930016e-05"
$ws.Range("O4").Value = [double]"8.459728460930016e-05"
$ws.Range("P4").Value = [double]"1.150973885098963"
$ws.Range("Q4").Value = "[0.672973801669885, 1.628973968528042]"
$ws.Range("R4").Value = [double]"1.514510869093932e-05"
$ws.Range("S4").Value = [double]"1.514510869093932e-05"
$ws.Range("T4").Value = [double]"55.87540383697127"
$ws.Range("U4").Value = "[43.40927954659462, 68.34152812734791]"
$ws.Range("V4").Value = [double]"1.158495521735858e-11"
$ws.Range("W4").Value = [double]"1.158495521735858e-11"
$ws.Range("X4").Value = [double]"20.33873873873911"
$ws.Range("Y4").Value = [double]"18.44444444444478"
$ws.Range("Z4").Value = [double]"22.23303303303344"

# Row 5
$ws.Range("B5").Value = [double]"1"
$ws.Range("H5").Value = [double]"4.951789011398144e-05"
$ws.Range("I5").Value = [double]"4.951789011398144e-05"
$ws.Range("L5").Value = [double]"40.38202916143509"
$ws.Range("M5").Value = "[18.52439236629607, 62.2396659565741]"
$ws.Range("N5").Value = [double]"0.0005489437908687833"
$ws.Range("O5").Value = [double]"0.0005489437908687833"
$ws.Range("P5").Value = [double]"1.050342288587577"
$ws.Range("Q5").Value = "[0.5220264069028087, 1.5786581702723463]"
$ws.Range("R5").Value = [double]"0.0002302196098986631"
$ws.Range("S5").Value = [double]"0.0002302196098986631"
$ws.Range("T5").Value = [double]"52.72373067356057"
$ws.Range("U5").Value = "[41.17014559655914, 64.27731575056201]"
$ws.Range("V5").Value = [double]"6.819433906457562e-12"
$ws.Range("W5").Value = [double]"6.819433906457562e-12"
$ws.Range("X5").Value = [double]"20.73753753753792"
$ws.Range("Y5").Value = [double]"18.64384384384419"
$ws.Range("Z5").Value = [double]"22.83123123123164"

# Row 6
$ws.Range("B6").Value = [double]"0"
$ws.Range("F6").Value = [double]"24.17000000000034"
$ws.Range("H6").Value = [double]"0.0007636207371958248"
$ws.Range("I6").Value = [double]"0.0007636207371958248"
$ws.Range("L6").Value = [double]"44.76317359729391"
$ws.Range("M6").Value = "[18.197412400639337, 71.32893479394848]"
$ws.Range("N6").Value = [double]"0.001447824523956198"
$ws.Range("O6").Value = [double]"0.001447824523956198"
$ws.Range("P6").Value = [double]"0.2075526678047313"
$ws.Range("Q6").Value = "[-0.5094474573388847, 0.9245527929483472]"
$ws.Range("R6").Value = [double]"0.5627837331447929"
$ws.Range("S6").Value = [double]"0.5627837331447929"
$ws.Range("T6").Value = [double]"66.16919687203156"
$ws.Range("U6").Value = "[50.720310387262245, 81.61808335680088]"
$ws.Range("V6").Value = [double]"4.298938982572054e-11"
$ws.Range("W6").Value = [double]"4.298938982572054e-11"
$ws.Range("X6").Value = [double]"23.37159159159192"
$ws.Range("Y6").Value = [double]"20.61345345345374"
$ws.Range("Z6").Value = [double]"26.12972972973009"

# Row 7
$ws.Range("F7").Value = [double]"24.17000000000034"
$ws.Range("H7").Value = [double]"0.00124550857186867"
$ws.Range("I7").Value = [double]"0.00124550857186867"
$ws.Range("L7").Value = [double]"40.27550311646542"
$ws.Range("M7").Value = "[13.216324942187981, 67.33468129074285]"
$ws.Range("N7").Value = [double]"0.004415519127883183"
$ws.Range("O7").Value = [double]"0.004415519127883183"
$ws.Range("P7").Value = [double]"0.4088158608275014"
$ws.Range("Q7").Value = "[-0.3081842643161137, 1.1258159859711165]"
$ws.Range("R7").Value = [double]"0.2568747515969052"
$ws.Range("S7").Value = [double]"0.2568747515969052"
$ws.Range("T7").Value = [double]"55.2709816179568"
$ws.Range("U7").Value = "[40.784804515799934, 69.75715872011367]"
$ws.Range("V7").Value = [double]"9.955991586707569e-10"
$ws.Range("W7").Value = [double]"9.955991586707569e-10"
$ws.Range("X7").Value = [double]"22.59737737737769"
$ws.Range("Y7").Value = [double]"19.83923923923952"
$ws.Range("Z7").Value = [double]"25.35551551551586"

# Row 8
$ws.Range("F8").Value = [double]"24.17000000000034"
$ws.Range("H8").Value = [double]"5.186642634413818e-05"
$ws.Range("I8").Value = [double]"5.186642634413818e-05"
$ws.Range("L8").Value = [double]"47.12467119680143"
$ws.Range("M8").Value = "[24.7734043601924, 69.47593803341046]"
$ws.Range("N8").Value = [double]"0.0001074181180742428"
$ws.Range("O8").Value = [double]"0.0001074181180742428"
$ws.Range("P8").Value = [double]"0.2327105669325773"
$ws.Range("Q8").Value = "[-0.34592111300788453, 0.8113422468730391]"
$ws.Range("R8").Value = [double]"0.4221917062175657"
$ws.Range("S8").Value = [double]"0.4221917062175657"
$ws.Range("T8").Value = [double]"57.89304651972073"
$ws.Range("U8").Value = "[44.46070487450227, 71.3253881649392]"
$ws.Range("V8").Value = [double]"3.5980329826657e-11"
$ws.Range("W8").Value = [double]"3.5980329826657e-11"
$ws.Range("X8").Value = [double]"23.27481481481514"
$ws.Range("Y8").Value = [double]"21.04894894894925"
$ws.Range("Z8").Value = [double]"25.50068068068104"

# Row 9
$ws.Range("F9").Value = [double]"24.17000000000034"
$ws.Range("H9").Value = [double]"2.308807435391635e-05"
$ws.Range("I9").Value = [double]"2.308807435391635e-05"
$ws.Range("L9").Value = [double]"51.00304546278252"
$ws.Range("M9").Value = "[27.97005788559487, 74.03603303997018]"
$ws.Range("N9").Value = [double]"5.418898009490469e-05"
$ws.Range("O9").Value = [double]"5.418898009490469e-05"
$ws.Range("P9").Value = [double]"0.1446579199851161"
$ws.Range("Q9").Value = "[-0.3710790121357297, 0.6603948521059619]"
$ws.Range("R9").Value = [double]"0.5749260782148493"
$ws.Range("S9").Value = [double]"0.5749260782148493"
$ws.Range("T9").Value = [double]"64.23599759561293"
$ws.Range("U9").Value = "[50.40296314249877, 78.06903204872708]"
$ws.Range("V9").Value = [double]"4.051869950671971e-12"
$ws.Range("W9").Value = [double]"4.051869950671971e-12"
$ws.Range("X9").Value = [double]"23.61353353353386"
$ws.Range("Y9").Value = [double]"21.62960960960992"
$ws.Range("Z9").Value = [double]"25.59745745745781"
